$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 27.50856438012401
$ws.Cells.Item(2, 5).Value = 26.6503791809082
$ws.Cells.Item(2, 6).Value = 29.63537164231805
$ws.Cells.Item(2, 7).Value = 25.18586788836676
$ws.Cells.Item(2, 8).Value = 60498713
$ws.Cells.Item(2, 9).Value = "UI"

$ws.Cells.Item(3, 4).Value = 30.14841529236411
$ws.Cells.Item(3, 5).Value = 30.05513572692871
$ws.Cells.Item(3, 6).Value = 33.58115596451929
$ws.Cells.Item(3, 7).Value = 28.51599798853331
$ws.Cells.Item(3, 8).Value = 60498713
$ws.Cells.Item(3, 9).Value = "UI"

$ws.Cells.Item(4, 4).Value = 31.61293091893109
$ws.Cells.Item(4, 5).Value = 27.21939659118652
$ws.Cells.Item(4, 6).Value = 33.54384737019176
$ws.Cells.Item(4, 7).Value = 26.77164820986108
$ws.Cells.Item(4, 8).Value = 60498713
$ws.Cells.Item(4, 9).Value = "UI"

$ws.Cells.Item(5, 4).Value = 29.77529084617096
$ws.Cells.Item(5, 5).Value = 27.6205005645752
$ws.Cells.Item(5, 6).Value = 29.79394747040486
$ws.Cells.Item(5, 7).Value = 24.01985393033624
$ws.Cells.Item(5, 8).Value = 60498713
$ws.Cells.Item(5, 9).Value = "UI"

$ws.Cells.Item(6, 4).Value = 30.88532995753701
$ws.Cells.Item(6, 5).Value = 33.22668075561523
$ws.Cells.Item(6, 6).Value = 34.13150730977915
$ws.Cells.Item(6, 7).Value = 29.90588047719969
$ws.Cells.Item(6, 8).Value = 60498713
$ws.Cells.Item(6, 9).Value = "UI"

$ws.Cells.Item(7, 4).Value = 36.03444386631087
$ws.Cells.Item(7, 5).Value = 41.71525573730469
$ws.Cells.Item(7, 6).Value = 41.90181485212737
$ws.Cells.Item(7, 7).Value = 34.56993275295333
$ws.Cells.Item(7, 8).Value = 60498713
$ws.Cells.Item(7, 9).Value = "UI"

$ws.Cells.Item(8, 4).Value = 49.90532904401199
$ws.Cells.Item(8, 5).Value = 48.90722274780273
$ws.Cells.Item(8, 6).Value = 51.30454387702167
$ws.Cells.Item(8, 7).Value = 45.00807712847811
$ws.Cells.Item(8, 8).Value = 60498713
$ws.Cells.Item(8, 9).Value = "UI"

$ws.Cells.Item(9, 4).Value = 54.65333317591352
$ws.Cells.Item(9, 5).Value = 58.20734024047852
$ws.Cells.Item(9, 6).Value = 58.76702476728819
$ws.Cells.Item(9, 7).Value = 52.7597301832071
$ws.Cells.Item(9, 8).Value = 60498713
$ws.Cells.Item(9, 9).Value = "UI"

$ws.Cells.Item(10, 4).Value = 46.91100543290828
$ws.Cells.Item(10, 5).Value = 48.0583610534668
$ws.Cells.Item(10, 6).Value = 48.97251422533325
$ws.Cells.Item(10, 7).Value = 44.96143281950716
$ws.Cells.Item(10, 8).Value = 60498713
$ws.Cells.Item(10, 9).Value = "UI"

$ws.Cells.Item(11, 4).Value = 48.95385851794768
$ws.Cells.Item(11, 5).Value = 50.838134765625
$ws.Cells.Item(11, 6).Value = 51.08066375111016
$ws.Cells.Item(11, 7).Value = 47.50800402069447
$ws.Cells.Item(11, 8).Value = 60498713
$ws.Cells.Item(11, 9).Value = "UI"

$ws.Cells.Item(12, 4).Value = 52.25601202054504
$ws.Cells.Item(12, 5).Value = 58.00212097167969
$ws.Cells.Item(12, 6).Value = 62.25573191934807
$ws.Cells.Item(12, 7).Value = 51.07134299796922
$ws.Cells.Item(12, 8).Value = 60498713
$ws.Cells.Item(12, 9).Value = "UI"

$ws.Cells.Item(13, 4).Value = 66.2574777497303
$ws.Cells.Item(13, 5).Value = 75.24976348876953
$ws.Cells.Item(13, 6).Value = 76.77024798983895
$ws.Cells.Item(13, 7).Value = 64.79296667504873
$ws.Cells.Item(13, 8).Value = 60498713
$ws.Cells.Item(13, 9).Value = "UI"

$ws.Cells.Item(14, 4).Value = 64.29858202245342
$ws.Cells.Item(14, 5).Value = 66.47203063964844
$ws.Cells.Item(14, 6).Value = 69.94208435434494
$ws.Cells.Item(14, 7).Value = 63.78553376799992
$ws.Cells.Item(14, 8).Value = 60498713
$ws.Cells.Item(14, 9).Value = "UI"

$ws.Cells.Item(15, 4).Value = 78.62652139563502
$ws.Cells.Item(15, 5).Value = 77.03141784667969
$ws.Cells.Item(15, 6).Value = 83.57973805380011
$ws.Cells.Item(15, 7).Value = 75.94003233162169
$ws.Cells.Item(15, 8).Value = 60498713
$ws.Cells.Item(15, 9).Value = "UI"

$ws.Cells.Item(16, 4).Value = 92.95059029885884
$ws.Cells.Item(16, 5).Value = 87.08504486083984
$ws.Cells.Item(16, 6).Value = 93.22188414064516
$ws.Cells.Item(16, 7).Value = 77.3839864709666
$ws.Cells.Item(16, 8).Value = 60498713
$ws.Cells.Item(16, 9).Value = "UI"

$ws.Cells.Item(17, 4).Value = 90.86158259766404
$ws.Cells.Item(17, 5).Value = 101.4563140869141
$ws.Cells.Item(17, 6).Value = 103.0970917521416
$ws.Cells.Item(17, 7).Value = 89.67084280052259
$ws.Cells.Item(17, 8).Value = 60498713
$ws.Cells.Item(17, 9).Value = "UI"

$ws.Cells.Item(18, 4).Value = 142.8464108152689
$ws.Cells.Item(18, 5).Value = 160.1115875244141
$ws.Cells.Item(18, 6).Value = 164.3386462323039
$ws.Cells.Item(18, 7).Value = 141.3622418156152
$ws.Cells.Item(18, 8).Value = 60498713
$ws.Cells.Item(18, 9).Value = "UI"

$ws.Cells.Item(19, 4).Value = 127.0537221880912
$ws.Cells.Item(19, 5).Value = 121.1527786254883
$ws.Cells.Item(19, 6).Value = 129.6794927946316
$ws.Cells.Item(19, 7).Value = 120.1081212933902
$ws.Cells.Item(19, 8).Value = 60498713
$ws.Cells.Item(19, 9).Value = "UI"

$ws.Cells.Item(20, 4).Value = 112.0277303954449
$ws.Cells.Item(20, 5).Value = 119.4641571044922
$ws.Cells.Item(20, 6).Value = 123.607040364981
$ws.Cells.Item(20, 7).Value = 104.7328615217077
$ws.Cells.Item(20, 8).Value = 60498713
$ws.Cells.Item(20, 9).Value = "UI"

$ws.Cells.Item(21, 4).Value = 179.8463372454641
$ws.Cells.Item(21, 5).Value = 154.4672698974609
$ws.Cells.Item(21, 6).Value = 182.3889717932381
$ws.Cells.Item(21, 7).Value = 153.9663063348354
$ws.Cells.Item(21, 8).Value = 60498713
$ws.Cells.Item(21, 9).Value = "UI"

$ws.Cells.Item(22, 4).Value = 132.8208672081967
$ws.Cells.Item(22, 5).Value = 153.4690551757812
$ws.Cells.Item(22, 6).Value = 158.6974070034868
$ws.Cells.Item(22, 7).Value = 125.0351742477011
$ws.Cells.Item(22, 8).Value = 60498713
$ws.Cells.Item(22, 9).Value = "UI"

$ws.Cells.Item(23, 4).Value = 164.7124382289049
$ws.Cells.Item(23, 5).Value = 175.8134460449219
$ws.Cells.Item(23, 6).Value = 180.4625814920453
$ws.Cells.Item(23, 7).Value = 161.4295688620132
$ws.Cells.Item(23, 8).Value = 60498713
$ws.Cells.Item(23, 9).Value = "UI"

$ws.Cells.Item(24, 4).Value = 158.6053664224804
$ws.Cells.Item(24, 5).Value = 176.4912567138672
$ws.Cells.Item(24, 6).Value = 197.0586115543464
$ws.Cells.Item(24, 7).Value = 155.6006090510449
$ws.Cells.Item(24, 8).Value = 60498713
$ws.Cells.Item(24, 9).Value = "UI"

$ws.Cells.Item(25, 4).Value = 264.4121339854667
$ws.Cells.Item(25, 5).Value = 293.3266906738281
$ws.Cells.Item(25, 6).Value = 308.1935189521645
$ws.Cells.Item(25, 7).Value = 224.8688831131373
$ws.Cells.Item(25, 8).Value = 60498713
$ws.Cells.Item(25, 9).Value = "UI"

$ws.Cells.Item(26, 4).Value = 284.4170971795236
$ws.Cells.Item(26, 5).Value = 272.0507202148437
$ws.Cells.Item(26, 6).Value = 286.0379907661139
$ws.Cells.Item(26, 7).Value = 260.1897049672482
$ws.Cells.Item(26, 8).Value = 60498713
$ws.Cells.Item(26, 9).Value = "UI"

$ws.Cells.Item(27, 4).Value = 300.0704551056804
$ws.Cells.Item(27, 5).Value = 298.9723510742188
$ws.Cells.Item(27, 6).Value = 307.3466218605492
$ws.Cells.Item(27, 7).Value = 275.0144620695706
$ws.Cells.Item(27, 8).Value = 60498713
$ws.Cells.Item(27, 9).Value = "UI"

$ws.Cells.Item(28, 4).Value = 287.187437119988
$ws.Cells.Item(28, 5).Value = 292.2766723632812
$ws.Cells.Item(28, 6).Value = 318.5550757236076
$ws.Cells.Item(28, 7).Value = 283.4853152352864
$ws.Cells.Item(28, 8).Value = 60498713
$ws.Cells.Item(28, 9).Value = "UI"

$ws.Cells.Item(29, 4).Value = 295.4949838652191
$ws.Cells.Item(29, 5).Value = 278.0114440917969
$ws.Cells.Item(29, 6).Value = 295.4949838652191
$ws.Cells.Item(29, 7).Value = 259.6076964604416
$ws.Cells.Item(29, 8).Value = 60498713
$ws.Cells.Item(29, 9).Value = "UI"

$ws.Cells.Item(30, 4).Value = 281.7207203420709
$ws.Cells.Item(30, 5).Value = 271.2086791992188
$ws.Cells.Item(30, 6).Value = 289.4365940813249
$ws.Cells.Item(30, 7).Value = 264.2903238373509
$ws.Cells.Item(30, 8).Value = 60498713
$ws.Cells.Item(30, 9).Value = "UI"

$ws.Cells.Item(31, 4).Value = 238.8291767061031
$ws.Cells.Item(31, 5).Value = 290.6048889160156
$ws.Cells.Item(31, 6).Value = 295.7207965700269
$ws.Cells.Item(31, 7).Value = 230.1774158778604
$ws.Cells.Item(31, 8).Value = 60498713
$ws.Cells.Item(31, 9).Value = "UI"

$ws.Cells.Item(32, 4).Value = 288.0441521366745
$ws.Cells.Item(32, 5).Value = 334.7273864746094
$ws.Cells.Item(32, 6).Value = 338.4728627368518
$ws.Cells.Item(32, 7).Value = 276.4022727459105
$ws.Cells.Item(32, 8).Value = 60498713
$ws.Cells.Item(32, 9).Value = "UI"

$ws.Cells.Item(33, 4).Value = 265.6674831396464
$ws.Cells.Item(33, 5).Value = 282.64599609375
$ws.Cells.Item(33, 6).Value = 287.3186971778324
$ws.Cells.Item(33, 7).Value = 255.2288739900677
$ws.Cells.Item(33, 8).Value = 60498713
$ws.Cells.Item(33, 9).Value = "UI"

$ws.Cells.Item(35, 4).Value = 172.1532735051223
$ws.Cells.Item(35, 5).Value = 172.8147583007812
$ws.Cells.Item(35, 6).Value = 184.1768625367101
$ws.Cells.Item(35, 7).Value = 166.3652221691902
$ws.Cells.Item(35, 8).Value = 60498713
$ws.Cells.Item(35, 9).Value = "UI"

$ws.Cells.Item(36, 4).Value = 141.4605322740552
$ws.Cells.Item(36, 5).Value = 118.5507965087891
$ws.Cells.Item(36, 6).Value = 141.763129708477
$ws.Cells.Item(36, 7).Value = 109.6094638370519
$ws.Cells.Item(36, 8).Value = 60498713
$ws.Cells.Item(36, 9).Value = "UI"

$ws.Cells.Item(37, 4).Value = 135.837086825381
$ws.Cells.Item(37, 5).Value = 123.4382858276367
$ws.Cells.Item(37, 6).Value = 137.2016418188173
$ws.Cells.Item(37, 7).Value = 120.758263570206
$ws.Cells.Item(37, 8).Value = 60498713
$ws.Cells.Item(37, 9).Value = "UI"

$ws.Cells.Item(38, 4).Value = 114.6686083276902
$ws.Cells.Item(38, 5).Value = 106.107421875
$ws.Cells.Item(38, 6).Value = 116.6510895646603
$ws.Cells.Item(38, 7).Value = 102.8131366548022
$ws.Cells.Item(38, 8).Value = 60498713
$ws.Cells.Item(38, 9).Value = "UI"

$ws.Cells.Item(39, 4).Value = 144.612356551601
$ws.Cells.Item(39, 5).Value = 183.8062744140625
$ws.Cells.Item(39, 6).Value = 187.2928080188765
$ws.Cells.Item(39, 7).Value = 143.9685404617738
$ws.Cells.Item(39, 8).Value = 60498713
$ws.Cells.Item(39, 9).Value = "UI"

$ws.Cells.Item(40, 4).Value = 219.5107683246197
$ws.Cells.Item(40, 5).Value = 263.9832458496094
$ws.Cells.Item(40, 6).Value = 268.225805032166
$ws.Cells.Item(40, 7).Value = 214.2348814711796
$ws.Cells.Item(40, 8).Value = 60498713
$ws.Cells.Item(40, 9).Value = "UI"

$ws.Cells.Item(41, 4).Value = 330.1005690317501
$ws.Cells.Item(41, 5).Value = 401.6668090820313
$ws.Cells.Item(41, 6).Value = 436.0862782064545
$ws.Cells.Item(41, 7).Value = 330.1005690317501
$ws.Cells.Item(41, 8).Value = 60498713
$ws.Cells.Item(41, 9).Value = "UI"

$ws.Cells.Item(42, 4).Value = 309.9521046378874
$ws.Cells.Item(42, 5).Value = 325.5070495605469
$ws.Cells.Item(42, 6).Value = 329.1564507405096
$ws.Cells.Item(42, 7).Value = 254.2634232934229
$ws.Cells.Item(42, 8).Value = 60498713
$ws.Cells.Item(42, 9).Value = "UI"

$ws.Cells.Item(43, 4).Value = 406.533254779058
$ws.Cells.Item(43, 5).Value = 434.8103637695313
$ws.Cells.Item(43, 6).Value = 473.0024237680694
$ws.Cells.Item(43, 7).Value = 392.714222870136
$ws.Cells.Item(43, 8).Value = 60498713
$ws.Cells.Item(43, 9).Value = "UI"

